$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "... тоталитарно [дистопично] общество ..." used to be split
# across three runs because "дистопично" was bracketed by a spell-check
# proofing mark (<w:proofErr type="spellStart"/> ... <w:proofErr
# type="spellEnd"/>). The fix removes that stale proofing split and merges
# the sentence back into a single contiguous run/text. A Find/Replace over
# the same literal text collapses the three runs (and drops the now
# pointless proofErr markers) into one run, matching the target XML.
# ---------------------------------------------------------------------------
$old1 = "владетел на тоталитарно дистопично общество, обитавано от"
$new1 = "владетел на тоталитарно дистопично общество, обитавано от"
$found1 = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)
if (-not $found1) {
    throw "Could not find the border-control paragraph text to normalize its runs."
}

Write-Output "Done."
